$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes all existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new listing's data
$ws.Range("A2").Value = "2024-08-06"
$ws.Range("B2").Value = "아이빔테크놀로지"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 223.4
$ws.Range("E2").Value = "삼성"
$ws.Range("F2").Value = 223.4
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2024-07-25"
$ws.Range("P2").Value = "2024-07-30"
$ws.Range("Q2").Value = 1675500

# Remove the old "그리드위즈" row, which is now at row 31 after the insert above
$ws.Rows.Item(31).Delete()
